# Weekly crime-data refresh for the 17th Precinct CompStat report.
# Bumps the report volume/number + the covered week dates, and refreshes
# every statistic cell in the Crime Complaints table (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 29   Number  44" -> "...45", and the week-covered line.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# ---------------------------------------------------------------------
# Helper: a few cells flip between a numeric value and the "N/A" text
# placeholders ("0" / "***.*"). Plain .Value assignment always coerces a
# digit-string back to a number and keeps the donor numeric style, so for
# those specific transitions we borrow formatting + value from a known
# same-type neighbour cell via Copy/PasteSpecial (xlPasteFormats=-4122,
# xlPasteValues=-4163), then for the reverse (text -> numeric) we just set
# the numeric style/value directly.
# ---------------------------------------------------------------------

function Set-AsTextPlaceholder($destAddr, $srcAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($destAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($destAddr).PasteSpecial(-4163) | Out-Null
}

# Row 15 (Rape): 1/2/-50% -> N/A (use row 23's stable placeholder cells)
Set-AsTextPlaceholder "C15" "C23"
Set-AsTextPlaceholder "D15" "D23"
Set-AsTextPlaceholder "E15" "E23"
$ws.Range("L15").Value = 27.272727272727

# Row 16 (Robbery)
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 79
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = 6.756756756756
$ws.Range("L16").Value = 43.636363636363
$ws.Range("M16").Value = 43.636363636363
$ws.Range("N16").Value = -85.714285714285

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = 4.040404040404
$ws.Range("L17").Value = 74.576271186440
$ws.Range("M17").Value = 110.204081632653
$ws.Range("N17").Value = -25.899280575539

# Row 18 (Burglary)
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 133
$ws.Range("K18").Value = 20.300751879699
$ws.Range("L18").Value = -0.621118012422
$ws.Range("M18").Value = 92.771084337349
$ws.Range("N18").Value = -84.962406015037

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 61
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 7.017543859649
$ws.Range("I19").Value = 577
$ws.Range("J19").Value = 468
$ws.Range("K19").Value = 23.290598290598
$ws.Range("L19").Value = 15.631262525050
$ws.Range("M19").Value = -10.681114551083
$ws.Range("N19").Value = -70.902672718103

# Row 20 (G.L.A.): C20 flips numeric -> N/A placeholder
Set-AsTextPlaceholder "C20" "C23"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 31.818181818181
$ws.Range("L20").Value = 41.463414634146
$ws.Range("N20").Value = -89.605734767025

# Row 21 (TOTAL)
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 22.222222222222
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -6.363636363636
$ws.Range("I21").Value = 991
$ws.Range("J21").Value = 826
$ws.Range("K21").Value = 19.975786924939
$ws.Range("L21").Value = 19.975786924939
$ws.Range("M21").Value = 14.170506912442
$ws.Range("N21").Value = -77.001624506846

# Row 22 (Transit): C/D/E/F flip from N/A placeholders -> numeric
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 18.181818181818
$ws.Range("L22").Value = 36.842105263157
$ws.Range("M22").Value = -13.333333333333

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -10.526315789473
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = -12.941176470588
$ws.Range("I24").Value = 1108
$ws.Range("J24").Value = 870
$ws.Range("K24").Value = 27.356321839080
$ws.Range("L24").Value = -3.231441048034
$ws.Range("M24").Value = 102.559414990859

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 64.705882352941
$ws.Range("I25").Value = 210
$ws.Range("J25").Value = 155
$ws.Range("K25").Value = 35.483870967741
$ws.Range("L25").Value = 47.887323943662
$ws.Range("M25").Value = 3.960396039603

# Row 26 (UCR Rape*): C/D/E flip numeric -> N/A placeholder
Set-AsTextPlaceholder "C26" "C23"
Set-AsTextPlaceholder "D26" "D23"
Set-AsTextPlaceholder "E26" "E23"
$ws.Range("L26").Value = -5.882352941176

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = -22.580645161290
$ws.Range("L27").Value = 20

# Row 30 (Hate Crimes): C30 flips numeric -> N/A placeholder
Set-AsTextPlaceholder "C30" "D30"
